# Updates to activity alignment
# Extends the time-series parameter tables (students, partnership,
# employment_smales, employment_sfemales, employment_couples) from
# 2027 out to 2035, and refreshes the projected "activity alignment"
# values for the three employment sheets over their respective
# recalculated ranges.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# students (sheet5) - flat extrapolation of the 2023 share (0.188)
# out to 2035
# ---------------------------------------------------------------
$wsStudents = $wb.Worksheets.Item("students")
for ($year = 2024; $year -le 2035; $year++) {
    $row = $year - 2010 + 2
    $wsStudents.Cells.Item($row, 1).Value = $year
    $wsStudents.Cells.Item($row, 2).Value = 0.188
}

# ---------------------------------------------------------------
# partnership (sheet6) - flat extrapolation of the 2023 share
# (0.6185) out to 2035
# ---------------------------------------------------------------
$wsPartnership = $wb.Worksheets.Item("partnership")
for ($year = 2024; $year -le 2035; $year++) {
    $row = $year - 2010 + 2
    $wsPartnership.Cells.Item($row, 1).Value = $year
    $wsPartnership.Cells.Item($row, 2).Value = 0.61850000000000005
}

# ---------------------------------------------------------------
# employment_smales (sheet7) - recalculated alignment from 2022
# onwards, extended out to 2035
# ---------------------------------------------------------------
$wsSMales = $wb.Worksheets.Item("employment_smales")
$smalesValues = @{
    2022 = 0.70015435151515204
    2023 = 0.70666020046620004
    2024 = 0.71316604941724904
    2025 = 0.71967189836829804
    2026 = 0.72617774731934703
    2027 = 0.73268359627039603
    2028 = 0.73918944522144503
    2029 = 0.74569529417249403
    2030 = 0.75220114312354303
    2031 = 0.75870699207459202
    2032 = 0.76521284102564102
    2033 = 0.77171868997669002
    2034 = 0.77822453892773902
    2035 = 0.78473038787878802
}
foreach ($year in ($smalesValues.Keys | Sort-Object)) {
    $row = $year - 2010 + 2
    $wsSMales.Cells.Item($row, 1).Value = $year
    $wsSMales.Cells.Item($row, 2).Value = $smalesValues[$year]
}

# ---------------------------------------------------------------
# employment_sfemales (sheet8) - recalculated alignment from 2023
# onwards, extended out to 2035
# ---------------------------------------------------------------
$wsSFemales = $wb.Worksheets.Item("employment_sfemales")
$sfemalesValues = @{
    2023 = 0.43197885000000003
    2024 = 0.4346719
    2025 = 0.43736494999999997
    2026 = 0.440058
    2027 = 0.44275104999999998
    2028 = 0.44544410000000001
    2029 = 0.44813714999999998
    2030 = 0.45083020000000001
    2031 = 0.45352324999999999
    2032 = 0.45621630000000002
    2033 = 0.45890934999999999
    2034 = 0.46160240000000002
    2035 = 0.46429545
}
foreach ($year in ($sfemalesValues.Keys | Sort-Object)) {
    $row = $year - 2010 + 2
    $wsSFemales.Cells.Item($row, 1).Value = $year
    $wsSFemales.Cells.Item($row, 2).Value = $sfemalesValues[$year]
}

# ---------------------------------------------------------------
# employment_couples (sheet9) - recalculated alignment from 2020
# onwards, extended out to 2035
# ---------------------------------------------------------------
$wsCouples = $wb.Worksheets.Item("employment_couples")
$couplesValues = @{
    2020 = 0.73668559333333306
    2021 = 0.73539384848484901
    2022 = 0.73410210363636397
    2023 = 0.73281035878787903
    2024 = 0.73151861393939399
    2025 = 0.73022686909090895
    2026 = 0.72893512424242402
    2027 = 0.72764337939393897
    2028 = 0.72635163454545504
    2029 = 0.72505988969696999
    2030 = 0.72376814484848495
    2031 = 0.72247640000000002
    2032 = 0.72247640000000002
    2033 = 0.72247640000000002
    2034 = 0.72247640000000002
    2035 = 0.72247640000000002
}
foreach ($year in ($couplesValues.Keys | Sort-Object)) {
    $row = $year - 2010 + 2
    $wsCouples.Cells.Item($row, 1).Value = $year
    $wsCouples.Cells.Item($row, 2).Value = $couplesValues[$year]
}

# ---------------------------------------------------------------
# Selections / active sheet - "students" becomes the active tab,
# each sheet's cursor/selection is refreshed to match the edited
# workbook.
# ---------------------------------------------------------------
$wsPartnership.Activate()
$wsPartnership.Range("A2").Select()

$wsSMales.Activate()
$wsSMales.Range("B2:B27").Select()

$wsSFemales.Activate()
$wsSFemales.Range("A2").Select()

$wsCouples.Activate()
$wsCouples.Range("A2").Select()

$wsStudents.Activate()
$wsStudents.Range("A2").Select()
